$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1762.4546
$ws.Range("I40").Value = 1237.4
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 1237.4
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -1062.4
$ws.Range("N40").Value = -2550

$ws.Range("H62").Value = 2881
$ws.Range("I62").Value = 2277.75
$ws.Range("J62").Value = 3182.625
$ws.Range("K62").Value = 2277.75
$ws.Range("L62").Value = 3182.625
$ws.Range("M62").Value = -1653.75
$ws.Range("N62").Value = -4430.625

$ws.Range("H65").Value = 2881
$ws.Range("I65").Value = 2277.75
$ws.Range("J65").Value = 3182.625
$ws.Range("K65").Value = 11388.75
$ws.Range("L65").Value = 15913.125
$ws.Range("M65").Value = -8268.75
$ws.Range("N65").Value = -22153.125

$ws.Range("H100").Value = 3763.95
$ws.Range("I100").Value = 3644.5386
$ws.Range("J100").Value = 3985.7144
$ws.Range("K100").Value = 3644.5386
$ws.Range("L100").Value = 3985.7144
$ws.Range("M100").Value = -3103.5386
$ws.Range("N100").Value = -5067.7144

$ws.Range("H103").Value = 830.7692
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 866.6667
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 2600.0001
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -3772.0001

$ws.Range("H106").Value = 2228.6316
$ws.Range("I106").Value = 1537
$ws.Range("J106").Value = 3414.2856
$ws.Range("K106").Value = 1537
$ws.Range("L106").Value = 3414.2856
$ws.Range("M106").Value = -906
$ws.Range("N106").Value = -4676.2856

$ws.Range("H132").Value = 7696357
$ws.Range("I132").Value = 9527795
$ws.Range("J132").Value = 4319.2
$ws.Range("K132").Value = 28583385
$ws.Range("L132").Value = 12957.6
$ws.Range("M132").Value = -28580855
$ws.Range("N132").Value = -18017.6

$ws.Range("H134").Value = 24421.6
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 24421.6
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 24421.6
$ws.Range("N134").Value = -34561.6

$ws.Range("H137").Value = 5126.1035
$ws.Range("I137").Value = 7092.533
$ws.Range("J137").Value = 3019.2144
$ws.Range("K137").Value = 21277.599
$ws.Range("L137").Value = 9057.643199999999
$ws.Range("M137").Value = -18727.599
$ws.Range("N137").Value = -14157.6432

$ws.Range("H138").Value = 4691.3706
$ws.Range("I138").Value = 2062.5293
$ws.Range("J138").Value = 6316.4727
$ws.Range("K138").Value = 6187.5879
$ws.Range("L138").Value = 18949.4181
$ws.Range("M138").Value = -1047.5879
$ws.Range("N138").Value = -29229.4181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14675.987
$ws.Range("I32").Value = 12724.086
$ws.Range("J32").Value = 20965.445
$ws.Range("K32").Value = 12724.086
$ws.Range("L32").Value = 20965.445
$ws.Range("M32").Value = -12437.086
$ws.Range("N32").Value = -21539.445

$ws.Range("H102").Value = 4053
$ws.Range("I102").Value = 2814.2727
$ws.Range("J102").Value = 6778.2
$ws.Range("K102").Value = 2814.2727
$ws.Range("L102").Value = 6778.2
$ws.Range("M102").Value = -1192.2727
$ws.Range("N102").Value = -10022.2

$ws.Range("H110").Value = 1686.96
$ws.Range("I110").Value = 638.625
$ws.Range("J110").Value = 3550.6667
$ws.Range("K110").Value = 638.625
$ws.Range("L110").Value = 3550.6667
$ws.Range("M110").Value = 1406.375
$ws.Range("N110").Value = -7640.6667

$ws.Range("H130").Value = 36886
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 36886
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 36886
$ws.Range("N130").Value = -46926

$ws.Range("H132").Value = 1972.1698
$ws.Range("I132").Value = 1489.7609
$ws.Range("J132").Value = 5142.2856
$ws.Range("K132").Value = 4469.2827
$ws.Range("L132").Value = 15426.8568
$ws.Range("M132").Value = -1939.2827
$ws.Range("N132").Value = -20486.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 29000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 29000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 29000
$ws.Range("N93").Value = -32744

$ws.Range("H94").Value = 1156.9048
$ws.Range("I94").Value = 653.2
$ws.Range("J94").Value = 2416.1667
$ws.Range("K94").Value = 653.2
$ws.Range("L94").Value = 2416.1667
$ws.Range("M94").Value = -202.2
$ws.Range("N94").Value = -3318.1667

$ws.Range("H107").Value = 2975
$ws.Range("I107").Value = 1462.5
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 1462.5
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 457.5
$ws.Range("N107").Value = -9840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10002320
$ws.Range("I134").Value = 13890839
$ws.Range("J134").Value = 3271.9285
$ws.Range("K134").Value = 41672517
$ws.Range("L134").Value = 9815.7855
$ws.Range("M134").Value = -41669982
$ws.Range("N134").Value = -14885.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 3072
$ws.Range("I57").Value = 499
$ws.Range("J57").Value = 3929.6667
$ws.Range("K57").Value = 1497
$ws.Range("L57").Value = 11789.0001
$ws.Range("M57").Value = -938
$ws.Range("N57").Value = -12907.0001

$ws.Range("H87").Value = 12494.6
$ws.Range("I87").Value = 5738
$ws.Range("J87").Value = 15390.286
$ws.Range("K87").Value = 17214
$ws.Range("L87").Value = 46170.858
$ws.Range("M87").Value = -15966
$ws.Range("N87").Value = -48666.858

$ws.Range("H90").Value = 12494.6
$ws.Range("I90").Value = 5738
$ws.Range("J90").Value = 15390.286
$ws.Range("K90").Value = 51642
$ws.Range("L90").Value = 138512.574
$ws.Range("M90").Value = -45402
$ws.Range("N90").Value = -150992.574

$ws.Range("H107").Value = 1860.6364
$ws.Range("I107").Value = 415.75
$ws.Range("J107").Value = 2686.2856
$ws.Range("K107").Value = 1247.25
$ws.Range("L107").Value = 8058.8568
$ws.Range("M107").Value = 672.75
$ws.Range("N107").Value = -11898.8568

$ws.Range("H113").Value = 1039.75
$ws.Range("I113").Value = 494.5
$ws.Range("J113").Value = 1100.3334
$ws.Range("K113").Value = 1483.5
$ws.Range("L113").Value = 3301.0002
$ws.Range("M113").Value = 686.5
$ws.Range("N113").Value = -7641.0002

$ws.Range("H121").Value = 11226.896
$ws.Range("I121").Value = 247.5
$ws.Range("J121").Value = 18977.059
$ws.Range("K121").Value = 742.5
$ws.Range("L121").Value = 56931.177
$ws.Range("M121").Value = 567.5
$ws.Range("N121").Value = -59551.177

$ws.Range("H136").Value = 2195.92
$ws.Range("I136").Value = 1444.875
$ws.Range("J136").Value = 3531.111
$ws.Range("K136").Value = 4334.625
$ws.Range("L136").Value = 10593.333
$ws.Range("M136").Value = 765.375
$ws.Range("N136").Value = -20793.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2925.926
$ws.Range("I80").Value = 3029.4119
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 3029.4119
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -2031.4119
$ws.Range("N80").Value = -4746

$ws.Range("H83").Value = 2925.926
$ws.Range("I83").Value = 3029.4119
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 15147.0595
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -10155.0595
$ws.Range("N83").Value = -23734

$ws.Range("H97").Value = 2083.1
$ws.Range("I97").Value = 1303.3334
$ws.Range("J97").Value = 3252.75
$ws.Range("K97").Value = 1303.3334
$ws.Range("L97").Value = 3252.75
$ws.Range("M97").Value = -807.3334
$ws.Range("N97").Value = -4244.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2267
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 4000.6667
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 4000.6667
$ws.Range("M22").Value = -238.3333
$ws.Range("N22").Value = -4590.6667

$ws.Range("H27").Value = 2267
$ws.Range("I27").Value = 533.3333
$ws.Range("J27").Value = 4000.6667
$ws.Range("K27").Value = 533.3333
$ws.Range("L27").Value = 4000.6667
$ws.Range("M27").Value = -426.3333
$ws.Range("N27").Value = -4214.6667

$ws.Range("H46").Value = 3200
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 5300
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 5300
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -5676

$ws.Range("H132").Value = 3006.0881
$ws.Range("I132").Value = 2606.6875
$ws.Range("J132").Value = 3361.111
$ws.Range("K132").Value = 7820.0625
$ws.Range("L132").Value = 10083.333
$ws.Range("M132").Value = -5290.0625
$ws.Range("N132").Value = -15143.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1711.7142
$ws.Range("I107").Value = 376.4
$ws.Range("J107").Value = 5050
$ws.Range("K107").Value = 1129.2
$ws.Range("L107").Value = 15150
$ws.Range("M107").Value = 790.8000000000002
$ws.Range("N107").Value = -18990

$ws.Range("H110").Value = 30580.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 30580.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 30580.5
$ws.Range("N110").Value = -38760.5

$ws.Range("H136").Value = 2490.6843
$ws.Range("I136").Value = 1963.7307
$ws.Range("J136").Value = 3632.4167
$ws.Range("K136").Value = 5891.1921
$ws.Range("L136").Value = 10897.2501
$ws.Range("M136").Value = -3341.1921
$ws.Range("N136").Value = -15997.2501

$ws.Range("H140").Value = 57450
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 57450
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 57450
$ws.Range("N140").Value = -67810
